$wb = $excel.ActiveWorkbook

# The "meta" sheet used to have an explicit "type" column (D) whose value was
# always "t" for every tab. The tab name is now used as the type instead, so
# that redundant column goes away and "title" shifts left into column D.
$ws = $wb.Worksheets.Item("meta")
$ws.Range("D1").EntireColumn.Delete()

# Make "meta" the active sheet/selection again (it had drifted to "set2").
$ws.Activate()
$ws.Range("A1").Select() | Out-Null
